$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "lxml" row (row 11) entirely; everything below shifts up.
$ws.Rows.Item(11).Delete()
